$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update inventory values
$ws.Range("D2").Value = 10
$ws.Range("D4").Value = 15

# Update the active cell selection
$ws.Range("G11").Select()
